# Scheduled data refresh: update market-price derived figures (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across the Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H63").Value = 64857.145
$ws.Range("I63").Value = 15000
$ws.Range("J63").Value = 68692.30499999999
$ws.Range("K63").Value = 15000
$ws.Range("L63").Value = 68692.30499999999
$ws.Range("M63").Value = -14376
$ws.Range("N63").Value = -69940.30499999999
$ws.Range("H66").Value = 64857.145
$ws.Range("I66").Value = 15000
$ws.Range("J66").Value = 68692.30499999999
$ws.Range("K66").Value = 45000
$ws.Range("L66").Value = 206076.915
$ws.Range("M66").Value = -41880
$ws.Range("N66").Value = -212316.915
$ws.Range("H68").Value = 58333.332
$ws.Range("J68").Value = 58333.332
$ws.Range("L68").Value = 58333.332
$ws.Range("N68").Value = -59831.332
$ws.Range("H71").Value = 58333.332
$ws.Range("J71").Value = 58333.332
$ws.Range("L71").Value = 174999.996
$ws.Range("N71").Value = -182487.996
$ws.Range("H129").Value = 1154.2858
$ws.Range("I129").Value = 1108.3334
$ws.Range("J129").Value = 1430
$ws.Range("K129").Value = 3325.0002
$ws.Range("L129").Value = 4290
$ws.Range("M129").Value = 1674.9998
$ws.Range("N129").Value = -14290
$ws.Range("H133").Value = 172998.83
$ws.Range("J133").Value = 172998.83
$ws.Range("L133").Value = 172998.83
$ws.Range("N133").Value = -183118.83
$ws.Range("H134").Value = 99400
$ws.Range("J134").Value = 99400
$ws.Range("L134").Value = 99400
$ws.Range("N134").Value = -109540
$ws.Range("H135").Value = 15625774
$ws.Range("I135").Value = 16129799
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 145168191
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -145165656
$ws.Range("N135").Value = -14070
$ws.Range("H137").Value = 1426
$ws.Range("I137").Value = 1337.2778
$ws.Range("J137").Value = 2224.5
$ws.Range("K137").Value = 4011.8334
$ws.Range("L137").Value = 6673.5
$ws.Range("M137").Value = -1461.8334
$ws.Range("N137").Value = -11773.5
$ws.Range("H138").Value = 1433.4359
$ws.Range("I138").Value = 927.88
$ws.Range("J138").Value = 2336.2144
$ws.Range("K138").Value = 2783.64
$ws.Range("L138").Value = 7008.6432
$ws.Range("M138").Value = 2356.36
$ws.Range("N138").Value = -17288.6432
$ws.Range("H140").Value = 74008.836
$ws.Range("J140").Value = 74008.836
$ws.Range("L140").Value = 74008.836
$ws.Range("N140").Value = -84368.836

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3011.5671
$ws.Range("I32").Value = 2637.4736
$ws.Range("K32").Value = 2637.4736
$ws.Range("M32").Value = -2350.4736
$ws.Range("H61").Value = 4486.25
$ws.Range("I61").Value = 3291
$ws.Range("K61").Value = 3291
$ws.Range("M61").Value = -3079
$ws.Range("H132").Value = 5265975.5
$ws.Range("I132").Value = 2054.4285
$ws.Range("J132").Value = 20004954
$ws.Range("K132").Value = 6163.2855
$ws.Range("L132").Value = 60014862
$ws.Range("M132").Value = -3633.2855
$ws.Range("N132").Value = -60019922
$ws.Range("H136").Value = 4486.25
$ws.Range("I136").Value = 3291
$ws.Range("K136").Value = 9873
$ws.Range("M136").Value = -7323

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 41667600
$ws.Range("I86").Value = 76923976
$ws.Range("J86").Value = 975.36365
$ws.Range("K86").Value = 76923976
$ws.Range("L86").Value = 975.36365
$ws.Range("M86").Value = -76922853
$ws.Range("N86").Value = -3221.36365
$ws.Range("H89").Value = 41667600
$ws.Range("I89").Value = 76923976
$ws.Range("J89").Value = 975.36365
$ws.Range("K89").Value = 384619880
$ws.Range("L89").Value = 4876.81825
$ws.Range("M89").Value = -384614264
$ws.Range("N89").Value = -16108.81825
$ws.Range("H94").Value = 4902906
$ws.Range("I94").Value = 6098330
$ws.Range("K94").Value = 6098330
$ws.Range("M94").Value = -6097879
$ws.Range("H107").Value = 42285.21
$ws.Range("I107").Value = 56098.43
$ws.Range("K107").Value = 56098.43
$ws.Range("M107").Value = -54178.43
$ws.Range("H138").Value = 86541
$ws.Range("J138").Value = 86541
$ws.Range("L138").Value = 86541
$ws.Range("N138").Value = -96821

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 105732.6
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 105732.6
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H75").Value = 57857.145
$ws.Range("J75").Value = 57857.145
$ws.Range("L75").Value = 57857.145
$ws.Range("N75").Value = -59853.145
$ws.Range("H78").Value = 57857.145
$ws.Range("J78").Value = 57857.145
$ws.Range("L78").Value = 173571.435
$ws.Range("N78").Value = -183555.435
$ws.Range("H99").Value = 3565.3333
$ws.Range("I99").Value = 3684.5334
$ws.Range("J99").Value = 2969.3333
$ws.Range("K99").Value = 3684.5334
$ws.Range("L99").Value = 2969.3333
$ws.Range("M99").Value = -2186.5334
$ws.Range("N99").Value = -5965.3333
$ws.Range("H126").Value = 3565.3333
$ws.Range("I126").Value = 3684.5334
$ws.Range("J126").Value = 2969.3333
$ws.Range("K126").Value = 11053.6002
$ws.Range("L126").Value = 8907.999899999999
$ws.Range("M126").Value = -8583.600199999999
$ws.Range("N126").Value = -13847.9999
$ws.Range("H135").Value = 66016.47
$ws.Range("J135").Value = 66016.47
$ws.Range("L135").Value = 66016.47
$ws.Range("N135").Value = -76156.47
$ws.Range("H138").Value = 66870.10000000001
$ws.Range("J138").Value = 66870.10000000001
$ws.Range("L138").Value = 66870.10000000001
$ws.Range("N138").Value = -77150.10000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 804.2941
$ws.Range("I14").Value = 804.2941
$ws.Range("K14").Value = 2412.8823
$ws.Range("M14").Value = -2239.8823
$ws.Range("H127").Value = 56559.125
$ws.Range("J127").Value = 56559.125
$ws.Range("L127").Value = 169677.375
$ws.Range("N127").Value = -179597.375
$ws.Range("H131").Value = 1694.5
$ws.Range("I131").Value = 1671
$ws.Range("K131").Value = 5013
$ws.Range("M131").Value = 27

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4250.913
$ws.Range("I40").Value = 4201.1333
$ws.Range("K40").Value = 4201.1333
$ws.Range("M40").Value = -4065.1333
$ws.Range("H55").Value = 1376.7858
$ws.Range("I55").Value = 446.8125
$ws.Range("K55").Value = 446.8125
$ws.Range("M55").Value = -273.8125
$ws.Range("H93").Value = 13891210
$ws.Range("I93").Value = 17546064
$ws.Range("J93").Value = 2768.6
$ws.Range("K93").Value = 17546064
$ws.Range("L93").Value = 2768.6
$ws.Range("M93").Value = -17544816
$ws.Range("N93").Value = -5264.6
$ws.Range("H136").Value = 9094320
$ws.Range("I136").Value = 3264.0557
$ws.Range("K136").Value = 9792.167099999999
$ws.Range("M136").Value = -7242.167099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 33010
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H107").Value = 11953
$ws.Range("J107").Value = 26901.428
$ws.Range("L107").Value = 80704.284
$ws.Range("N107").Value = -84544.284
$ws.Range("H136").Value = 1713.4546
$ws.Range("I136").Value = 1307.8334
$ws.Range("J136").Value = 2200.2
$ws.Range("K136").Value = 3923.5002
$ws.Range("L136").Value = 6600.599999999999
$ws.Range("M136").Value = -1373.5002
$ws.Range("N136").Value = -11700.6
$ws.Range("H139").Value = 70747
$ws.Range("J139").Value = 70747
$ws.Range("L139").Value = 70747
$ws.Range("N139").Value = -81027
$ws.Range("H141").Value = 90765
$ws.Range("J141").Value = 90765
$ws.Range("L141").Value = 90765
$ws.Range("N141").Value = -101125
